$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.351.91'
$ws.Range('E2').Value = '  -3.37%  '
$ws.Range('D3').Value = '3.164.57'
$ws.Range('E3').Value = '  -2.57%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.43'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.46'
$ws.Range('E6').Value = '  -5.69%  '
$ws.Range('D8').Value = '3.162.45'
$ws.Range('E8').Value = '  -2.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('E10').Value = '  -5.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.48'
$ws.Range('E11').Value = '  -7.78%  '
$ws.Range('E12').Value = '  -5.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  -6.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.64'
$ws.Range('E14').Value = '  -9.03%  '
$ws.Range('D15').Value = '3.680.55'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('D16').Value = '64.361.62'
$ws.Range('E16').Value = '  -3.40%  '
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '3.164.33'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.95'
$ws.Range('E19').Value = '  -5.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '482.12'
$ws.Range('E20').Value = '  -4.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.73'
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('E22').Value = '  -4.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.77'
$ws.Range('E23').Value = '  -3.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.76'
$ws.Range('E24').Value = '  -6.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.67'
$ws.Range('E25').Value = '  -3.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.87'
$ws.Range('E27').Value = '  -4.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.49'
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  -7.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.74'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.113'
$ws.Range('E31').Value = '  -19.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.75'
$ws.Range('E32').Value = '  -3.45%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.30'
$ws.Range('E34').Value = '  -6.23%  '
$ws.Range('E35').Value = '  -3.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.83'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.00'
$ws.Range('E37').Value = '  -6.35%  '
$ws.Range('D38').Value = '0.0₃0729'
$ws.Range('E38').Value = '  -7.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '451.85'
$ws.Range('E39').Value = '  -8.25%  '
$ws.Range('E40').Value = '  -10.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0397'
$ws.Range('E41').Value = '  -6.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.45'
$ws.Range('E42').Value = '  -4.10%  '
$ws.Range('E43').Value = '  -8.54%  '
$ws.Range('D44').Value = '2.845.82'
$ws.Range('E44').Value = '  -3.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.267'
$ws.Range('E45').Value = '  -8.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.28'
$ws.Range('E46').Value = '  -7.85%  '
$ws.Range('E47').Value = '  -7.15%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  -7.88%  '
$ws.Range('E50').Value = '  -4.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.29'
$ws.Range('E51').Value = '  -1.67%  '
